$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 367/368, pushing all existing data (old rows 367-460)
# down to rows 369-462. This also extends the used range to A1:R462.
$ws.Range("A367:R368").Insert()

# Fill in the new weekly record (week of 2021-11-08, serial date 44508)
# for "1a (cosecha)" / "2a (cosecha)" from "Región de Arica y Parinacota".
$ws.Range("A367").Value = 8
$ws.Range("B367").Value = "Terminal La Palmera de La Serena"
$ws.Range("C367").Value = "Coquimbo"
$ws.Range("D367").Value = 44508
$ws.Range("E367").Value = 4
$ws.Range("F367").Value = 100112004
$ws.Range("G367").Value = "Cebolla"
$ws.Range("H367").Value = "Sin especificar"
$ws.Range("I367").Value = "1a (cosecha)"
$ws.Range("J367").Value = 3360
$ws.Range("K367").Value = 4800
$ws.Range("L367").Value = 5000
$ws.Range("M367").Value = 4900
$ws.Range("N367").Value = "`$/malla 18 kilos"
$ws.Range("O367").Value = "Región de Arica y Parinacota"
$ws.Range("P367").Value = 272
$ws.Range("Q367").Value = 18
$ws.Range("R367").Value = "Hortaliza"

$ws.Range("A368").Value = 8
$ws.Range("B368").Value = "Terminal La Palmera de La Serena"
$ws.Range("C368").Value = "Coquimbo"
$ws.Range("D368").Value = 44508
$ws.Range("E368").Value = 4
$ws.Range("F368").Value = 100112004
$ws.Range("G368").Value = "Cebolla"
$ws.Range("H368").Value = "Sin especificar"
$ws.Range("I368").Value = "2a (cosecha)"
$ws.Range("J368").Value = 1600
$ws.Range("K368").Value = 4500
$ws.Range("L368").Value = 4600
$ws.Range("M368").Value = 4550
$ws.Range("N368").Value = "`$/malla 18 kilos"
$ws.Range("O368").Value = "Región de Arica y Parinacota"
$ws.Range("P368").Value = 253
$ws.Range("Q368").Value = 18
$ws.Range("R368").Value = "Hortaliza"

Write-Output "Done"
